$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Quantity ("F" column, "Stock") updates -------------------------------
$ws.Range("F11").Value  = 61
$ws.Range("F36").Value  = 157
$ws.Range("F41").Value  = 68
$ws.Range("F44").Value  = 42
$ws.Range("F46").Value  = 7
$ws.Range("F49").Value  = 21
$ws.Range("F58").Value  = 36
$ws.Range("F67").Value  = 46
$ws.Range("F77").Value  = 10
$ws.Range("F89").Value  = 92
$ws.Range("F94").Value  = 60
$ws.Range("F100").Value = 38
$ws.Range("F120").Value = 23

# --- Price ("D" column) updates -------------------------------------------
$ws.Range("D33").Value  = 158.75
$ws.Range("D128").Value = 7.67
$ws.Range("D129").Value = 7.67

# --- F129 text changes from "15" to "11" (kept as text, matches F128 type) -
$ws.Range("F129").Value = "11"

# --- Number formatting: the whole "F" (stock) column is now shown as a
# plain integer ("0") instead of the previous thousand-separated / text
# formats, while keeping each cell's existing alignment untouched. ---------
$ws.Range("F128:F129").NumberFormat = "0"
$ws.Range("F2:F127").NumberFormat = "0"
$ws.Range("F1").NumberFormat      = "0"

# --- Restore the selection left behind on the sheet on save ---------------
[void]$ws.Range("G3").Select()
